# Fruta / hortaliza, semanal
# Insert a new weekly price record as a new row 264 (pushing the existing
# rows 264-315 down to 265-316), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 264, shifting rows 264:315 -> 265:316
$ws.Rows.Item(264).Insert()

# Populate the newly inserted row 264 with the new weekly record
$ws.Range("A264").Value2 = 10
$ws.Range("B264").Value2 = "Vega Modelo de Temuco"
$ws.Range("C264").Value2 = "La Araucanía"
$ws.Range("D264").Value2 = 45015
$ws.Range("E264").Value2 = 9
$ws.Range("F264").Value2 = "Fruta"
$ws.Range("G264").Value2 = 100101
$ws.Range("H264").Value2 = "Berries"
$ws.Range("I264").Value2 = 100112025
$ws.Range("J264").Value2 = "Frutilla"
$ws.Range("K264").Value2 = "Sin especificar"
$ws.Range("L264").Value2 = "Primera"
$ws.Range("M264").Value2 = 35
$ws.Range("N264").Value2 = 7000
$ws.Range("O264").Value2 = 7000
$ws.Range("P264").Value2 = 7000
$ws.Range("Q264").Value2 = "$/caja 7 kilos"
$ws.Range("R264").Value2 = "Región de La Araucanía"
$ws.Range("S264").Value2 = 1000
$ws.Range("T264").Value2 = 7
